$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added to the "Jengibre" (Hortaliza,
# Vega Modelo de Temuco) dataset. It belongs right after the existing
# row 71 (chronologically it is inserted as the new row 72), pushing all
# the subsequent rows (old 72..109) down by one (new 73..110).
$ws.Rows("72").Insert()

# Populate the newly inserted row 72 with the new record.
$ws.Cells.Item(72, 1).Value = 10
$ws.Cells.Item(72, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(72, 3).Value = "La Araucanía"
$ws.Cells.Item(72, 4).Value = 44455
$ws.Cells.Item(72, 5).Value = 9
$ws.Cells.Item(72, 6).Value = 100114007
$ws.Cells.Item(72, 7).Value = "Jengibre"
$ws.Cells.Item(72, 8).Value = "Sin especificar"
$ws.Cells.Item(72, 9).Value = "Primera"
$ws.Cells.Item(72, 10).Value = 50
$ws.Cells.Item(72, 11).Value = 20000
$ws.Cells.Item(72, 12).Value = 20000
$ws.Cells.Item(72, 13).Value = 20000
$ws.Cells.Item(72, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(72, 15).Value = "Perú"
$ws.Cells.Item(72, 16).Value = 1538
$ws.Cells.Item(72, 17).Value = 13
$ws.Cells.Item(72, 18).Value = "Hortaliza"
